$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts), recomputed to replace the old "Strike#" values.
# Write the new values for rows 2-16.
$newK = @{
    2  = 1
    3  = 4
    4  = 8
    5  = 4
    6  = 9
    7  = 2
    8  = 6
    9  = 2
    10 = 5
    11 = 6
    12 = 3
    13 = 7
    14 = 3
    15 = 3
    16 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
